$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# Fill in Protein (Q) and Ontogeny (R) for existing row 2 (TestPopulation)
$ws.Range("Q2").Value = "CYP3A4, CYP2D6"
$ws.Range("R2").Value = "CYP3A4, CYP2D6"

# Add new row 3 (TestPopulation_noOnto), mirroring row 2 but without Protein/Ontogeny
$ws.Range("A3").Value = "TestPopulation_noOnto"
$ws.Range("B3").Value = "Human"
$ws.Range("C3").Value = "European_ICRP_2002"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0
$ws.Range("H3").Value = "kg"
$ws.Range("K3").Value = "cm"
$ws.Range("L3").Value = 22
$ws.Range("M3").Value = 41
$ws.Range("P3").Value = "kg/m²"

$ws.Range("R3").Select()
